# Update TPM-derived NATMI metrics for the Col3a1-Mag LR-pair sheet.
# Only numeric result columns (G:T, excluding the rank columns K/L which
# are unaffected) change; identifiers in columns A-F/K/L stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value  = 22.628972
$ws.Range("H2").Value  = 67.886916
$ws.Range("I2").Value  = 0.004372730881336598
$ws.Range("J2").Value  = 0.004372730881336598
$ws.Range("M2").Value  = 0.232947
$ws.Range("N2").Value  = 0.698841
$ws.Range("O2").Value  = 0.2572219815457369
$ws.Range("P2").Value  = 0.2572219815457369
$ws.Range("Q2").Value  = 5.271351140484001
$ws.Range("R2").Value  = 47.442160264356
$ws.Range("S2").Value  = 0.001124762502063636
$ws.Range("T2").Value  = 0.001124762502063636

# Row 3 (ECs -> Resolving-Mac)
$ws.Range("G3").Value  = 22.628972
$ws.Range("H3").Value  = 67.886916
$ws.Range("I3").Value  = 0.004372730881336598
$ws.Range("J3").Value  = 0.004372730881336598
$ws.Range("M3").Value  = 0.6726793333333333
$ws.Range("N3").Value  = 2.018038
$ws.Range("O3").Value  = 0.7427780184542632
$ws.Range("P3").Value  = 0.7427780184542632
$ws.Range("Q3").Value  = 15.22204179897867
$ws.Range("R3").Value  = 136.998376190808
$ws.Range("S3").Value  = 0.003247968379272962
$ws.Range("T3").Value  = 0.003247968379272962

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value  = 0.9598063873258337
$ws.Range("J4").Value  = 0.9598063873258338
$ws.Range("M4").Value  = 0.232947
$ws.Range("N4").Value  = 0.698841
$ws.Range("O4").Value  = 0.2572219815457369
$ws.Range("P4").Value  = 0.2572219815457369
$ws.Range("Q4").Value  = 1157.051881712727
$ws.Range("R4").Value  = 10413.46693541454
$ws.Range("S4").Value  = 0.246883300848206
$ws.Range("T4").Value  = 0.246883300848206

# Row 5 (FAPs -> Resolving-Mac)
$ws.Range("I5").Value  = 0.9598063873258337
$ws.Range("J5").Value  = 0.9598063873258338
$ws.Range("M5").Value  = 0.6726793333333333
$ws.Range("N5").Value  = 2.018038
$ws.Range("O5").Value  = 0.7427780184542632
$ws.Range("P5").Value  = 0.7427780184542632
$ws.Range("Q5").Value  = 3341.210182670719
$ws.Range("R5").Value  = 30070.89164403647
$ws.Range("S5").Value  = 0.7129230864776277
$ws.Range("T5").Value  = 0.7129230864776278

# Row 6 (MuSCs -> ECs)
$ws.Range("G6").Value  = 182.6322073333333
$ws.Range("H6").Value  = 547.896622
$ws.Range("I6").Value  = 0.0352911079183418
$ws.Range("J6").Value  = 0.0352911079183418
$ws.Range("M6").Value  = 0.232947
$ws.Range("N6").Value  = 0.698841
$ws.Range("O6").Value  = 0.2572219815457369
$ws.Range("P6").Value  = 0.2572219815457369
$ws.Range("Q6").Value  = 42.543624801678
$ws.Range("R6").Value  = 382.892623215102
$ws.Range("S6").Value  = 0.009077648709700323
$ws.Range("T6").Value  = 0.009077648709700323

# Row 7 (MuSCs -> Resolving-Mac)
$ws.Range("G7").Value  = 182.6322073333333
$ws.Range("H7").Value  = 547.896622
$ws.Range("I7").Value  = 0.0352911079183418
$ws.Range("J7").Value  = 0.0352911079183418
$ws.Range("M7").Value  = 0.6726793333333333
$ws.Range("N7").Value  = 2.018038
$ws.Range("O7").Value  = 0.7427780184542632
$ws.Range("P7").Value  = 0.7427780184542632
$ws.Range("Q7").Value  = 122.8529114741818
$ws.Range("R7").Value  = 1105.676203267636
$ws.Range("S7").Value  = 0.02621345920864148
$ws.Range("T7").Value  = 0.02621345920864148

# Row 8 (Resolving-Mac -> ECs)
$ws.Range("G8").Value  = 2.741590666666667
$ws.Range("H8").Value  = 8.224772
$ws.Range("I8").Value  = 0.00052977387448787
$ws.Range("J8").Value  = 0.00052977387448787
$ws.Range("M8").Value  = 0.232947
$ws.Range("N8").Value  = 0.698841
$ws.Range("O8").Value  = 0.2572219815457369
$ws.Range("P8").Value  = 0.2572219815457369
$ws.Range("Q8").Value  = 0.6386453210280001
$ws.Range("R8").Value  = 5.747807889252
$ws.Range("S8").Value  = 0.0001362694857669324
$ws.Range("T8").Value  = 0.0001362694857669324

# Row 9 (Resolving-Mac -> Resolving-Mac)
$ws.Range("G9").Value  = 2.741590666666667
$ws.Range("H9").Value  = 8.224772
$ws.Range("I9").Value  = 0.00052977387448787
$ws.Range("J9").Value  = 0.00052977387448787
$ws.Range("M9").Value  = 0.6726793333333333
$ws.Range("N9").Value  = 2.018038
$ws.Range("O9").Value  = 0.7427780184542632
$ws.Range("P9").Value  = 0.7427780184542632
$ws.Range("Q9").Value  = 1.844211381926222
$ws.Range("R9").Value  = 16.597902437336
$ws.Range("S9").Value  = 0.0003935043887209376
$ws.Range("T9").Value  = 0.0003935043887209376
